$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# New row 54: 四方坪站 data for 2025-11-27
$ws.Cells.Item(54, 1).Value = 45988
$ws.Cells.Item(54, 2).Value = "四方坪站"
$ws.Cells.Item(54, 3).Value = 8224.7999999999993
$ws.Cells.Item(54, 4).Value = 7348.22
$ws.Cells.Item(54, 5).Value = 2762.38
$ws.Cells.Item(54, 6).Value = 341

# New row 55: 高岭站 data for 2025-11-27
$ws.Cells.Item(55, 1).Value = 45988
$ws.Cells.Item(55, 2).Value = "高岭站"
$ws.Cells.Item(55, 3).Value = 4637.62
$ws.Cells.Item(55, 4).Value = 3909.14
$ws.Cells.Item(55, 5).Value = 1192.3800000000001
$ws.Cells.Item(55, 6).Value = 169

# Scroll / selection state, matching the saved view in the workbook
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("I56").Select()
